$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 55
$ws.Range("A55").Value = 130838040
$ws.Range("B55").Value = 83228
$ws.Range("E55").Value = 1467
$ws.Range("F55").Value = "Rödbrun blekspik"
$ws.Range("G55").Value = "Sclerophora coniophaea"
$ws.Range("H55").Value = "(Norman) J.Mattsson & Middelb."
$ws.Range("Q55").Value = 445709
$ws.Range("R55").Value = 7026357
$ws.Range("Z55").Value = "11:11"
$ws.Range("AB55").Value = "11:11"
$ws.Range("AC55").Value = "På död gren i hålighet vid basen av gammal levande grov gran (42 cm dbh) i gammal granskog"

# Row 56
$ws.Range("A56").Value = 130838768
$ws.Range("B56").Value = 79243
$ws.Range("E56").Value = 6425
$ws.Range("F56").Value = "Garnlav"
$ws.Range("G56").Value = "Alectoria sarmentosa"
$ws.Range("H56").Value = "(Ach.) Ach."
$ws.Range("Q56").Value = 445697
$ws.Range("R56").Value = 7026283
$ws.Range("Z56").Value = "11:56"
$ws.Range("AB56").Value = "11:56"
$ws.Range("AC56").Value = "På gammal gran i gammal barrblandskog"

# Row 76
$ws.Range("A76").Value = 130838833
$ws.Range("B76").Value = 89193
$ws.Range("E76").Value = 510
$ws.Range("F76").Value = "Doftskinn"
$ws.Range("G76").Value = "Cystostereum murrayi"
$ws.Range("H76").Value = "(Berk. & M.A.Curtis.) Pouzar"
$ws.Range("Q76").Value = 445685
$ws.Range("R76").Value = 7026259
$ws.Range("S76").Value = 4
$ws.Range("Z76").Value = "12:07"
$ws.Range("AB76").Value = "12:07"
$ws.Range("AC76").Value = "På granlåga i gammal granskog"

# Row 77
$ws.Range("A77").Value = 130837733
$ws.Range("B77").Value = 79243
$ws.Range("D77").Value = "NT"
$ws.Range("E77").Value = 6425
$ws.Range("F77").Value = "Garnlav"
$ws.Range("G77").Value = "Alectoria sarmentosa"
$ws.Range("H77").Value = "(Ach.) Ach."
$ws.Range("Q77").Value = 445720
$ws.Range("R77").Value = 7026343
$ws.Range("S77").Value = 10
$ws.Range("Z77").Value = "10:59"
$ws.Range("AB77").Value = "10:59"
$ws.Range("AC77").Value = "På gammal död gran i gammal granskog"

# Row 78
$ws.Range("A78").Value = 130837541
$ws.Range("B78").Value = 75221
$ws.Range("D78").Value = "LC"
$ws.Range("E78").Value = 6428
$ws.Range("F78").Value = "Rostfläck"
$ws.Range("G78").Value = "Arthonia vinosa"
$ws.Range("H78").Value = "Leight."
$ws.Range("Q78").Value = 445740
$ws.Range("R78").Value = 7026322
$ws.Range("S78").Value = 8
$ws.Range("Z78").Value = "10:52"
$ws.Range("AB78").Value = "10:52"
$ws.Range("AC78").Value = "På tunna kvistar vid basen på gammal levande gran"
